{"js": "// Replace each two-digit multiplication equation with its updated result.\n// The mapping below is derived exactly from the authoritative diff: each\n// old equation string is unique in the document, so a direct text search\n// + replace per pair is unambiguous.\nconst replacements = [\n  [\"86\u00d788=7568\", \"46\u00d730=1380\"],\n  [\"91\u00d752=4732\", \"42\u00d760=2520\"],\n  [\"21\u00d784=1764\", \"32\u00d761=1952\"],\n  [\"52\u00d726=1352\", \"88\u00d729=2552\"],\n  [\"17\u00d711=187\", \"83\u00d733=2739\"],\n  [\"68\u00d786=5848\", \"24\u00d762=1488\"],\n  [\"57\u00d734=1938\", \"76\u00d731=2356\"],\n  [\"36\u00d752=1872\", \"45\u00d743=1935\"],\n  [\"62\u00d777=4774\", \"59\u00d751=3009\"],\n  [\"64\u00d720=1280\", \"35\u00d725=875\"],\n  [\"26\u00d715=390\", \"84\u00d783=6972\"],\n  [\"76\u00d752=3952\", \"50\u00d747=2350\"],\n  [\"80\u00d766=5280\", \"25\u00d770=1750\"],\n  [\"87\u00d736=3132\", \"67\u00d745=3015\"],\n  [\"36\u00d776=2736\", \"66\u00d755=3630\"],\n  [\"96\u00d763=6048\", \"12\u00d781=972\"],\n  [\"88\u00d771=6248\", \"31\u00d753=1643\"],\n  [\"61\u00d718=1098\", \"68\u00d755=3740\"],\n  [\"13\u00d716=208\", \"61\u00d749=2989\"],\n  [\"49\u00d753=2597\", \"47\u00d781=3807\"],\n  [\"50\u00d774=3700\", \"12\u00d785=1020\"],\n  [\"14\u00d783=1162\", \"74\u00d776=5624\"],\n  [\"72\u00d788=6336\", \"81\u00d770=5670\"],\n  [\"63\u00d740=2520\", \"13\u00d722=286\"],\n  [\"86\u00d792=7912\", \"96\u00d714=1344\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update each two-digit multiplication equation to its new result.\n# Mapping is derived exactly from the target diff; every \"Old\" string is\n# unique within the document body, so Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @{ Old = \"86\u00d788=7568\"; New = \"46\u00d730=1380\" },\n  @{ Old = \"91\u00d752=4732\"; New = \"42\u00d760=2520\" },\n  @{ Old = \"21\u00d784=1764\"; New = \"32\u00d761=1952\" },\n  @{ Old = \"52\u00d726=1352\"; New = \"88\u00d729=2552\" },\n  @{ Old = \"17\u00d711=187\"; New = \"83\u00d733=2739\" },\n  @{ Old = \"68\u00d786=5848\"; New = \"24\u00d762=1488\" },\n  @{ Old = \"57\u00d734=1938\"; New = \"76\u00d731=2356\" },\n  @{ Old = \"36\u00d752=1872\"; New = \"45\u00d743=1935\" },\n  @{ Old = \"62\u00d777=4774\"; New = \"59\u00d751=3009\" },\n  @{ Old = \"64\u00d720=1280\"; New = \"35\u00d725=875\" },\n  @{ Old = \"26\u00d715=390\"; New = \"84\u00d783=6972\" },\n  @{ Old = \"76\u00d752=3952\"; New = \"50\u00d747=2350\" },\n  @{ Old = \"80\u00d766=5280\"; New = \"25\u00d770=1750\" },\n  @{ Old = \"87\u00d736=3132\"; New = \"67\u00d745=3015\" },\n  @{ Old = \"36\u00d776=2736\"; New = \"66\u00d755=3630\" },\n  @{ Old = \"96\u00d763=6048\"; New = \"12\u00d781=972\" },\n  @{ Old = \"88\u00d771=6248\"; New = \"31\u00d753=1643\" },\n  @{ Old = \"61\u00d718=1098\"; New = \"68\u00d755=3740\" },\n  @{ Old = \"13\u00d716=208\"; New = \"61\u00d749=2989\" },\n  @{ Old = \"49\u00d753=2597\"; New = \"47\u00d781=3807\" },\n  @{ Old = \"50\u00d774=3700\"; New = \"12\u00d785=1020\" },\n  @{ Old = \"14\u00d783=1162\"; New = \"74\u00d776=5624\" },\n  @{ Old = \"72\u00d788=6336\"; New = \"81\u00d770=5670\" },\n  @{ Old = \"63\u00d740=2520\"; New = \"13\u00d722=286\" },\n  @{ Old = \"86\u00d792=7912\"; New = \"96\u00d714=1344\" }\n)\n\nforeach ($pair in $replacements) {\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $pair.Old\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $pair.New\n\n  $found = $find.Execute(\n    $pair.Old,   # FindText\n    $false,      # MatchCase\n    $false,      # MatchWholeWord\n    $false,      # MatchWildcards\n    $false,      # MatchSoundsLike\n    $false,      # MatchAllWordForms\n    $true,       # Forward\n    1,           # Wrap (wdFindContinue)\n    $false,      # Format\n    $pair.New,   # ReplaceWith\n    2            # Replace (wdReplaceAll)\n  )\n\n  if (-not $found) {\n    throw \"No match found for: $($pair.Old)\"\n  }\n}\n"}
